$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# Insert 8 new rows before the current row 5 (rows 5-12 new, shifting old row5.. down by 8)
$ws.Rows("5:12").Insert()

# Update header row 2: D2 changes from "Tue" to "Est"
$ws.Range("D2").Value = "Est"

# Row 3: remove D3 value (it used to be "Hours", now blank)
$ws.Range("D3").ClearContents()

# New rows 5-11 content
$ws.Range("A5").Value = "Complete tests for Sask CrownRoyalty Rate"
$ws.Range("B5").Value = "A"
$ws.Range("D5").Value = 5

$ws.Range("A6").Value = "Brake out Used Royalty Rate"
$ws.Range("B6").Value = "A"
$ws.Range("D6").Value = 5

$ws.Range("A7").Value = "Refactor DB Load with test"
$ws.Range("B7").Value = "K"
$ws.Range("D7").Value = 14

$ws.Range("A8").Value = "Write the our version of orm"
$ws.Range("B8").Value = "K"
$ws.Range("D8").Value = 14

$ws.Range("A9").Value = "Read organization of tests"
$ws.Range("B9").Value = "K"
$ws.Range("D9").Value = 1

$ws.Range("A10").Value = "Complete tests IORG1995"
$ws.Range("B10").Value = "A"
$ws.Range("D10").Value = 10

$ws.Range("A11").Value = "Finish IOGR1995"
$ws.Range("B11").Value = "A"
$ws.Range("D11").Value = 5

# Apply wrap text formatting to the new column A cells, matching existing rows (style index 2)
$ws.Range("A5:A11").WrapText = $true

# Row 33 (previously row 25): add "Done" status in column C
$ws.Range("C33").Value = "Done"

# Row 34 (previously row 26): replace content entirely
$ws.Range("A34").Value = "Create Data Access Strategy Step 1 example"
$ws.Range("B34").Value = "L"
$ws.Range("C34").Value = "Dane"
$ws.Range("D34").Value = 10

# Update selection to D12
$ws.Range("D12").Select()
